$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The conversation has ended without a clear decision on which movie to show on Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding which movie to show on Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not finalized.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded, and it indicates that no movie was selected in this meeting.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The conversation has ended without making a decision about which movie will be shown on Friday. Therefore, no action will be taken regarding the acquisition of movie rights.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected in this meeting.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision-making process ended without a clear choice for Friday’s movie, so I have recorded a no-decision outcome.`n"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The rights to `"Barbie`" have been successfully acquired.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision process concluded without a clear choice of a movie for Friday, so no action will be taken regarding acquiring movie rights.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded as a no decision about the movie to be shown on Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be selected for Friday as there was no consensus reached.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("D23").Value = "both_movies, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: It appears that we did not reach a decision regarding which movie to show on Friday, leading to the conclusion that no decision can be made at this time.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for the movie `"Barbie`" has been successfully made.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision to acquire the movie `"Barbie`" has been recorded successfully.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to show on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision made is that there was no agreement on a movie to show on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no decision made.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been made.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Oppenheimer.`"`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision to acquire the rights to show both movies has been recorded.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded and the rights to `"Barbie`" have been acquired.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday's showing.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been registered as no movie was selected for Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision process did not result in a specific movie being selected for Friday. Therefore, I will proceed with the no_decision function.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to both movies have been acquired.`n"
$ws.Range("D45").Value = "both_movies, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie has concluded with no selection made.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be shown on Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Oppenheimer`" will be shown on Friday.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie for Friday.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday ended without a definitive choice.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie has resulted in no agreement.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday, so I have recorded that as no decision.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no decision.`"`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made, and thus no action will be taken to acquire movie rights.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies, `"Oppenheimer`" and `"Barbie,`" has been recorded successfully.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been made that there will be no selection for the movie to be shown on Friday.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision about Friday's movie has not been made, and therefore, no acquisition will occur.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The rights to both movies, `"Barbie`" and `"Oppenheimer,`" have been successfully acquired.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be reached, so there will be no movie selected.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie was made.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Oppenheimer.`"`n"

Write-Host "Updated cells"